$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextValue "D2" "66.448.38"
Set-TextValue "E2" "  -1.76%  "
Set-TextValue "D3" "3.437.60"
Set-TextValue "E3" "  -1.55%  "
Set-TextValue "E4" "  +0.04%  "
Set-TextValue "D5" "580.25"
Set-TextValue "E5" "  -2.73%  "
Set-TextValue "D6" "174.49"
Set-TextValue "E6" "  -2.43%  "
Set-TextValue "E7" "  -0.05%  "
Set-TextValue "E8" "  +1.25%  "
Set-TextValue "D9" "3.437.24"
Set-TextValue "E9" "  -1.54%  "
Set-TextValue "E10" "  -3.71%  "
Set-TextValue "E11" "  -3.37%  "
Set-TextValue "E12" "  -2.81%  "
Set-TextValue "D13" "4.029.62"
Set-TextValue "E13" "  -1.68%  "
Set-TextValue "D14" "30.87"
Set-TextValue "E14" "  -3.80%  "
Set-TextValue "E15" "  -3.63%  "
Set-TextValue "D16" "66.469.85"
Set-TextValue "E16" "  -1.70%  "
Set-TextValue "E17" "  -3.33%  "
Set-TextValue "D18" "3.439.93"
Set-TextValue "E18" "  -1.65%  "
Set-TextValue "E19" "  -4.26%  "
Set-TextValue "D20" "13.79"
Set-TextValue "E20" "  -3.79%  "
Set-TextValue "D21" "373.91"
Set-TextValue "E21" "  -5.14%  "
Set-TextValue "E22" "  -1.71%  "
Set-TextValue "D23" "1.00"
Set-TextValue "E23" "  +0.08%  "
Set-TextValue "E24" "  -0.03%  "
Set-TextValue "D25" "70.82"
Set-TextValue "E25" "  -3.10%  "
Set-TextValue "D26" "0.525"
Set-TextValue "E26" "  -1.95%  "
Set-TextValue "E27" "  -2.57%  "
Set-TextValue "D28" "9.81"
Set-TextValue "E28" "  -6.34%  "
Set-TextValue "E29" "  -2.48%  "
Set-TextValue "E30" "  +0.25%  "
Set-TextValue "D31" "5.87"
Set-TextValue "E31" "  -4.86%  "
Set-TextValue "E32" "  -3.14%  "
Set-TextValue "E33" "  +0.49%  "
Set-TextValue "E34" "  -6.27%  "
Set-TextValue "D35" "0.999"
Set-TextValue "E36" "  -3.54%  "
Set-TextValue "E37" "  -5.44%  "
Set-TextValue "D38" "160.18"
Set-TextValue "E38" "  -2.17%  "
Set-TextValue "E39" "  +0.09%  "
Set-TextValue "D40" "27.22"
Set-TextValue "E40" "  +4.55%  "
Set-TextValue "E41" "  -5.63%  "
Set-TextValue "D42" "2.62"
Set-TextValue "E42" "  -3.98%  "
Set-TextValue "D43" "6.55"
Set-TextValue "E43" "  -5.25%  "
Set-TextValue "E44" "  -4.93%  "
Set-TextValue "D45" "2.681.94"
Set-TextValue "E45" "  -5.90%  "
Set-TextValue "D46" "0.0691"
Set-TextValue "E46" "  -4.81%  "
Set-TextValue "D47" "25.26"
Set-TextValue "E47" "  -4.70%  "
Set-TextValue "D48" "40.58"
Set-TextValue "E48" "  -3.48%  "
Set-TextValue "E49" "  -2.90%  "
Set-TextValue "D50" "317.85"
Set-TextValue "E50" "  -5.54%  "
Set-TextValue "E51" "  -4.62%  "
